$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B slightly (stored width goes from 10.7109375 to 12 characters).
# The host's ColumnWidth -> stored-width conversion adds a small fixed offset,
# so 11.14 round-trips to exactly 12 in the saved XML.
$ws.Columns.Item(2).ColumnWidth = 11.14

# Zoom the view to 85% and move the active selection to I28
$excel.ActiveWindow.Zoom = 85
$ws.Range("I28").Select()
